$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data: week 4 / problem 14890 "경사로"
$ws.Range("A16").Value = 14890
$ws.Range("B16").Value = "경사로"
$ws.Range("C16").Value = 45910

# Reuse C15's existing cell format (built-in date format already in the
# style table) rather than letting NumberFormat mint a brand-new custom
# numFmt entry.
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null

# Move the active selection to reflect the newly added row, matching the
# author's saved view state (cursor parked one row below the new data).
$ws.Range("A17").Select() | Out-Null
